# Append the new row (row 4) of tech/quiz round scores to the "Scores" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Gagan Naik"
$ws.Range("B4").Value = "mentalmilestone24@gmail.com"
$ws.Range("C4").Value = 34
